# Daily GitHub-Actions refresh of the cryptos price/volume snapshot.
# Updates Price (D) / Volume(1h) (E) figures for most rows, and for the
# WrappedEther/Polkadot pair (rows 13-14) the ranking also swapped places
# so the coin name + link + price + volume all move together.
# Numeric-looking Price strings ("4.21", "65.51", ...) are written with a
# leading apostrophe so Excel keeps them as text (matching the workbook's
# original inline-string cells) instead of auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.783.72'
$ws.Range("E2").Value = '  +0.08%  '

$ws.Range("D3").Value = '1.648.86'
$ws.Range("E3").Value = '  -0.06%  '

$ws.Range("E4").Value = '  +0.79%  '

$ws.Range("E5").Value = '  +0.91%  '

$ws.Range("E6").Value = '  +0.18%  '

$ws.Range("E7").Value = '  +0.79%  '

$ws.Range("E8").Value = '  +0.10%  '

$ws.Range("E9").Value = '  -0.24%  '

$ws.Range("E10").Value = '  +0.11%  '

$ws.Range("E11").Value = '  +0.32%  '

$ws.Range("D12").Value = '1.873.83'
$ws.Range("E12").Value = '  -0.09%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = "'4.21"
$ws.Range("E13").Value = '  +0.95%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.659.73'
$ws.Range("E14").Value = '  +0.86%  '

$ws.Range("D15").Value = "'0.532"
$ws.Range("E15").Value = '  -0.22%  '

$ws.Range("D16").Value = "'65.51"
$ws.Range("E16").Value = '  -0.57%  '

$ws.Range("D17").Value = '26.795.05'
$ws.Range("E17").Value = '  +0.13%  '

$ws.Range("E18").Value = '  -0.43%  '

$ws.Range("D19").Value = "'217.51"
$ws.Range("E19").Value = '  -0.91%  '

$ws.Range("E20").Value = '  +0.73%  '

$ws.Range("D21").Value = "'2.54"
$ws.Range("E21").Value = '  +18.50%  '

$ws.Range("E22").Value = '  +0.28%  '

$ws.Range("D23").Value = "'6.32"
$ws.Range("E23").Value = '  -0.30%  '

$ws.Range("E24").Value = '  +0.17%  '

$ws.Range("D25").Value = "'145.67"
$ws.Range("E25").Value = '  -1.26%  '

$ws.Range("E26").Value = '  +0.83%  '

$ws.Range("D27").Value = "'0.121"
$ws.Range("E27").Value = '  -0.22%  '

$ws.Range("E28").Value = '  +3.61%  '

$ws.Range("E29").Value = '  -0.15%  '

$ws.Range("E30").Value = '  +0.25%  '

$ws.Range("E31").Value = '  +0.99%  '

$ws.Range("E32").Value = '  -0.78%  '

$ws.Range("E33").Value = '  -0.14%  '

$ws.Range("E34").Value = '  +1.50%  '

$ws.Range("D35").Value = '1.278.18'
$ws.Range("E35").Value = '  +0.59%  '

$ws.Range("E36").Value = '  +1.90%  '

$ws.Range("E37").Value = '  +1.30%  '

$ws.Range("E38").Value = '  +5.35%  '

$ws.Range("D39").Value = "'0.834"
$ws.Range("E39").Value = '  +2.45%  '

$ws.Range("E40").Value = '  +0.77%  '

$ws.Range("E41").Value = '  +2.05%  '

$ws.Range("E42").Value = '  -1.23%  '

$ws.Range("D43").Value = "'5.41"
$ws.Range("E43").Value = '  +0.77%  '

$ws.Range("D44").Value = '1.799.28'
$ws.Range("E44").Value = '  +0.92%  '

$ws.Range("D45").Value = "'92.12"

$ws.Range("D46").Value = "'59.86"
$ws.Range("E46").Value = '  +7.52%  '

$ws.Range("E47").Value = '  +0.94%  '

$ws.Range("E48").Value = '  +2.45%  '

$ws.Range("E49").Value = '  +0.36%  '

$ws.Range("D50").Value = "'7.81"
$ws.Range("E50").Value = '  +1.86%  '

$ws.Range("E51").Value = '  +1.24%  '
